# Remove the leftover "setor" placeholder rectangle from slide 2.
# The shape (id=3, name="Rectangle 4") held a stray "setor" label that
# shouldn't have shipped with the slide; delete the whole shape.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "Rectangle 4" -and $shape.HasTextFrame -and $shape.TextFrame.TextRange.Text.Trim() -eq "setor") {
        $shape.Delete()
        break
    }
}
